$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 20000
$ws.Range("I21").Value = 20000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 20000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -19532
$ws.Range("H23").Value = 20000
$ws.Range("I23").Value = 20000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -19766
$ws.Range("H31").Value = 83335450
$ws.Range("I31").Value = 83335450
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 250006350
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -250006120
$ws.Range("H87").Value = 74000
$ws.Range("I87").Value = 74000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 74000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -72752
$ws.Range("H90").Value = 74000
$ws.Range("I90").Value = 74000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 222000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -215760
$ws.Range("H141").Value = 4030.8333
$ws.Range("I141").Value = 4063.3333
$ws.Range("J141").Value = 3998.3333
$ws.Range("K141").Value = 12189.9999
$ws.Range("L141").Value = 11994.9999
$ws.Range("M141").Value = -7009.999899999999
$ws.Range("N141").Value = -22354.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 131
$ws.Range("I5").Value = 131
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 131
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -19
$ws.Range("N5").ClearContents()
$ws.Range("H74").Value = 787.25
$ws.Range("I74").Value = 716.3333
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 716.3333
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = 157.6667
$ws.Range("N74").Value = -2748
$ws.Range("H77").Value = 787.25
$ws.Range("I77").Value = 716.3333
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 3581.6665
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = 786.3334999999997
$ws.Range("N77").Value = -13736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 131
$ws.Range("I4").Value = 131
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 131
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -16
$ws.Range("N4").ClearContents()
$ws.Range("H33").Value = 26500
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 26500
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 26500
$ws.Range("N33").Value = -27172
$ws.Range("M33").ClearContents()
$ws.Range("H64").Value = 830.44446
$ws.Range("I64").Value = 946.6667
$ws.Range("J64").Value = 772.3333
$ws.Range("K64").Value = 946.6667
$ws.Range("L64").Value = 772.3333
$ws.Range("M64").Value = -721.6667
$ws.Range("N64").Value = -1222.3333
$ws.Range("H67").Value = 830.44446
$ws.Range("I67").Value = 946.6667
$ws.Range("J67").Value = 772.3333
$ws.Range("K67").Value = 946.6667
$ws.Range("L67").Value = 772.3333
$ws.Range("M67").Value = -166.6667
$ws.Range("N67").Value = -2332.3333
$ws.Range("H107").Value = 3419.6
$ws.Range("I107").Value = 3419.6
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3419.6
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1499.6
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 129999.5
$ws.Range("I23").Value = 129999.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 129999.5
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -129759.5
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 129999.5
$ws.Range("I27").Value = 129999.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 129999.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -129807.5
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 998
$ws.Range("I31").Value = 997.3333
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 997.3333
$ws.Range("L31").Value = 1000
$ws.Range("M31").Value = -702.3333
$ws.Range("N31").Value = -1590
$ws.Range("H34").Value = 998
$ws.Range("I34").Value = 997.3333
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 997.3333
$ws.Range("L34").Value = 1000
$ws.Range("M34").Value = -795.3333
$ws.Range("N34").Value = -1404
$ws.Range("H99").Value = 3727.75
$ws.Range("I99").Value = 2708.25
$ws.Range("J99").Value = 4747.25
$ws.Range("K99").Value = 2708.25
$ws.Range("L99").Value = 4747.25
$ws.Range("M99").Value = -1210.25
$ws.Range("N99").Value = -7743.25
$ws.Range("H126").Value = 3727.75
$ws.Range("I126").Value = 2708.25
$ws.Range("J126").Value = 4747.25
$ws.Range("K126").Value = 8124.75
$ws.Range("L126").Value = 14241.75
$ws.Range("M126").Value = -5654.75
$ws.Range("N126").Value = -19181.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 373.5
$ws.Range("I34").Value = 197
$ws.Range("J34").Value = 550
$ws.Range("K34").Value = 591
$ws.Range("L34").Value = 1650
$ws.Range("M34").Value = -507
$ws.Range("N34").Value = -1818
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H103").Value = 843.8889
$ws.Range("I103").Value = 39.25
$ws.Range("J103").Value = 1487.6
$ws.Range("K103").Value = 117.75
$ws.Range("L103").Value = 4462.799999999999
$ws.Range("M103").Value = 761.25
$ws.Range("N103").Value = -6220.799999999999
$ws.Range("H109").Value = 1063.8572
$ws.Range("I109").Value = 507.83334
$ws.Range("J109").Value = 4400
$ws.Range("K109").Value = 1523.50002
$ws.Range("L109").Value = 13200
$ws.Range("M109").Value = -483.5000199999999
$ws.Range("N109").Value = -15280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 21030.834
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 21030.834
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 21030.834
$ws.Range("N95").Value = -26522.834
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3530
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4638
$ws.Range("I132").Value = 4327.6665
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 12982.9995
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -10452.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H34").Value = 15000
$ws.Range("I34").Value = 15000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -14828
$ws.Range("N34").ClearContents()
$ws.Range("H93").Value = 3196.5715
$ws.Range("I93").Value = 2854.3333
$ws.Range("J93").Value = 5250
$ws.Range("K93").Value = 2854.3333
$ws.Range("L93").Value = 5250
$ws.Range("M93").Value = -1606.3333
$ws.Range("H110").Value = 150000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 150000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 150000
$ws.Range("N110").Value = -158180
$ws.Range("H136").Value = 3498.4
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 3998
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 11994
$ws.Range("M136").Value = -1950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 693.7778
$ws.Range("I107").Value = 575
$ws.Range("J107").Value = 931.3333
$ws.Range("K107").Value = 1725
$ws.Range("L107").Value = 2793.9999
$ws.Range("M107").Value = 195
$ws.Range("N107").Value = -6633.9999
$ws.Range("H132").Value = 2653.111
$ws.Range("I132").Value = 2672.25
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 8016.75
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -5486.75
$ws.Range("N132").Value = -12560
